$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The three BHRUT_ITU (col I) figures for 48/49/50 were provisional estimates
#     (shown with a yellow highlight); the real reported figure turned out to be 0,
#     so clear the highlight and correct the values. ---
$ws.Range("I54:I56").ClearFormats()
$ws.Range("I54").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("I56").Value = 0

# --- Append the newly-reported daily rows 57:76 ---
$ws.Range("A57").Value = 43951
$ws.Range("B57").Value = 144
$ws.Range("C57").Value = 39
$ws.Range("D57").Formula = "=B57-C57"
$ws.Range("E57").Value = 259
$ws.Range("F57").Value = 110
$ws.Range("G57").Formula = "=E57-F57"
$ws.Range("H57").Value = 36
$ws.Range("I57").Value = 2
$ws.Range("J57").Formula = "=H57-I57"
$ws.Range("K57").Formula = "=B57+E57+H57"
$ws.Range("L57").Formula = "=C57+F57+I57"
$ws.Range("M57").Formula = "=K57-L57"

$ws.Range("A58").Value = 43952
$ws.Range("B58").Value = 138
$ws.Range("C58").Value = 40
$ws.Range("D58").Formula = "=B58-C58"
$ws.Range("E58").Value = 242
$ws.Range("F58").Value = 105
$ws.Range("G58").Formula = "=E58-F58"
$ws.Range("H58").Value = 32
$ws.Range("I58").Value = 2
$ws.Range("J58").Formula = "=H58-I58"
$ws.Range("K58").Formula = "=B58+E58+H58"
$ws.Range("L58").Formula = "=C58+F58+I58"
$ws.Range("M58").Formula = "=K58-L58"

$ws.Range("A59").Value = 43953
$ws.Range("B59").Value = 134
$ws.Range("C59").Value = 37
$ws.Range("D59").Formula = "=B59-C59"
$ws.Range("E59").Value = 237
$ws.Range("F59").Value = 89
$ws.Range("G59").Formula = "=E59-F59"
$ws.Range("H59").Value = 30
$ws.Range("I59").Value = 2
$ws.Range("J59").Formula = "=H59-I59"
$ws.Range("K59").Formula = "=B59+E59+H59"
$ws.Range("L59").Formula = "=C59+F59+I59"
$ws.Range("M59").Formula = "=K59-L59"

$ws.Range("A60").Value = 43954
$ws.Range("B60").Value = 124
$ws.Range("C60").Value = 37
$ws.Range("D60").Formula = "=B60-C60"
$ws.Range("E60").Value = 219
$ws.Range("F60").Value = 80
$ws.Range("G60").Formula = "=E60-F60"
$ws.Range("H60").Value = 31
$ws.Range("I60").Value = 2
$ws.Range("J60").Formula = "=H60-I60"
$ws.Range("K60").Formula = "=B60+E60+H60"
$ws.Range("L60").Formula = "=C60+F60+I60"
$ws.Range("M60").Formula = "=K60-L60"

$ws.Range("A61").Value = 43955
$ws.Range("B61").Value = 123
$ws.Range("C61").Value = 35
$ws.Range("D61").Formula = "=B61-C61"
$ws.Range("E61").Value = 209
$ws.Range("F61").Value = 73
$ws.Range("G61").Formula = "=E61-F61"
$ws.Range("H61").Value = 32
$ws.Range("I61").Value = 2
$ws.Range("J61").Formula = "=H61-I61"
$ws.Range("K61").Formula = "=B61+E61+H61"
$ws.Range("L61").Formula = "=C61+F61+I61"
$ws.Range("M61").Formula = "=K61-L61"

$ws.Range("A62").Value = 43956
$ws.Range("B62").Value = 120
$ws.Range("C62").Value = 33
$ws.Range("D62").Formula = "=B62-C62"
$ws.Range("E62").Value = 203
$ws.Range("F62").Value = 79
$ws.Range("G62").Formula = "=E62-F62"
$ws.Range("H62").Value = 25
$ws.Range("I62").Value = 12
$ws.Range("J62").Formula = "=H62-I62"
$ws.Range("K62").Formula = "=B62+E62+H62"
$ws.Range("L62").Formula = "=C62+F62+I62"
$ws.Range("M62").Formula = "=K62-L62"

$ws.Range("A63").Value = 43957
$ws.Range("B63").Value = 104
$ws.Range("C63").Value = 30
$ws.Range("D63").Formula = "=B63-C63"
$ws.Range("E63").Value = 192
$ws.Range("F63").Value = 71
$ws.Range("G63").Formula = "=E63-F63"
$ws.Range("H63").Value = 23
$ws.Range("I63").Value = 10
$ws.Range("J63").Formula = "=H63-I63"
$ws.Range("K63").Formula = "=B63+E63+H63"
$ws.Range("L63").Formula = "=C63+F63+I63"
$ws.Range("M63").Formula = "=K63-L63"

$ws.Range("A64").Value = 43958
$ws.Range("B64").Value = 92
$ws.Range("C64").Value = 29
$ws.Range("D64").Formula = "=B64-C64"
$ws.Range("E64").Value = 176
$ws.Range("F64").Value = 70
$ws.Range("G64").Formula = "=E64-F64"
$ws.Range("H64").Value = 24
$ws.Range("I64").Value = 11
$ws.Range("J64").Formula = "=H64-I64"
$ws.Range("K64").Formula = "=B64+E64+H64"
$ws.Range("L64").Formula = "=C64+F64+I64"
$ws.Range("M64").Formula = "=K64-L64"

$ws.Range("A65").Value = 43959
$ws.Range("B65").Value = 82
$ws.Range("C65").Value = 27
$ws.Range("D65").Formula = "=B65-C65"
$ws.Range("E65").Value = 176
$ws.Range("F65").Value = 59
$ws.Range("G65").Formula = "=E65-F65"
$ws.Range("H65").Value = 24
$ws.Range("I65").Value = 12
$ws.Range("J65").Formula = "=H65-I65"
$ws.Range("K65").Formula = "=B65+E65+H65"
$ws.Range("L65").Formula = "=C65+F65+I65"
$ws.Range("M65").Formula = "=K65-L65"

$ws.Range("A66").Value = 43960
$ws.Range("B66").Value = 75
$ws.Range("C66").Value = 28
$ws.Range("D66").Formula = "=B66-C66"
$ws.Range("E66").Value = 173
$ws.Range("F66").Value = 58
$ws.Range("G66").Formula = "=E66-F66"
$ws.Range("H66").Value = 25
$ws.Range("I66").Value = 11
$ws.Range("J66").Formula = "=H66-I66"
$ws.Range("K66").Formula = "=B66+E66+H66"
$ws.Range("L66").Formula = "=C66+F66+I66"
$ws.Range("M66").Formula = "=K66-L66"

$ws.Range("A67").Value = 43961
$ws.Range("B67").Value = 72
$ws.Range("C67").Value = 24
$ws.Range("D67").Formula = "=B67-C67"
$ws.Range("E67").Value = 175
$ws.Range("F67").Value = 55
$ws.Range("G67").Formula = "=E67-F67"
$ws.Range("H67").Value = 23
$ws.Range("I67").Value = 11
$ws.Range("J67").Formula = "=H67-I67"
$ws.Range("K67").Formula = "=B67+E67+H67"
$ws.Range("L67").Formula = "=C67+F67+I67"
$ws.Range("M67").Formula = "=K67-L67"

$ws.Range("A68").Value = 43962
$ws.Range("B68").Value = 71
$ws.Range("C68").Value = 26
$ws.Range("D68").Formula = "=B68-C68"
$ws.Range("E68").Value = 168
$ws.Range("F68").Value = 46
$ws.Range("G68").Formula = "=E68-F68"
$ws.Range("H68").Value = 21
$ws.Range("I68").Value = 10
$ws.Range("J68").Formula = "=H68-I68"
$ws.Range("K68").Formula = "=B68+E68+H68"
$ws.Range("L68").Formula = "=C68+F68+I68"
$ws.Range("M68").Formula = "=K68-L68"

$ws.Range("A69").Value = 43963
$ws.Range("B69").Value = 74
$ws.Range("C69").Value = 27
$ws.Range("D69").Formula = "=B69-C69"
$ws.Range("E69").Value = 171
$ws.Range("F69").Value = 36
$ws.Range("G69").Formula = "=E69-F69"
$ws.Range("H69").Value = 21
$ws.Range("I69").Value = 8
$ws.Range("J69").Formula = "=H69-I69"
$ws.Range("K69").Formula = "=B69+E69+H69"
$ws.Range("L69").Formula = "=C69+F69+I69"
$ws.Range("M69").Formula = "=K69-L69"

$ws.Range("A70").Value = 43964
$ws.Range("B70").Value = 72
$ws.Range("C70").Value = 26
$ws.Range("D70").Formula = "=B70-C70"
$ws.Range("E70").Value = 162
$ws.Range("F70").Value = 32
$ws.Range("G70").Formula = "=E70-F70"
$ws.Range("H70").Value = 21
$ws.Range("I70").Value = 8
$ws.Range("J70").Formula = "=H70-I70"
$ws.Range("K70").Formula = "=B70+E70+H70"
$ws.Range("L70").Formula = "=C70+F70+I70"
$ws.Range("M70").Formula = "=K70-L70"

$ws.Range("A71").Value = 43965
$ws.Range("B71").Value = 69
$ws.Range("C71").Value = 23
$ws.Range("D71").Formula = "=B71-C71"
$ws.Range("E71").Value = 152
$ws.Range("F71").Value = 33
$ws.Range("G71").Formula = "=E71-F71"
$ws.Range("H71").Value = 19
$ws.Range("I71").Value = 8
$ws.Range("J71").Formula = "=H71-I71"
$ws.Range("K71").Formula = "=B71+E71+H71"
$ws.Range("L71").Formula = "=C71+F71+I71"
$ws.Range("M71").Formula = "=K71-L71"

$ws.Range("A72").Value = 43966
$ws.Range("B72").Value = 67
$ws.Range("C72").Value = 24
$ws.Range("D72").Formula = "=B72-C72"
$ws.Range("E72").Value = 161
$ws.Range("F72").Value = 40
$ws.Range("G72").Formula = "=E72-F72"
$ws.Range("H72").Value = 21
$ws.Range("I72").Value = 8
$ws.Range("J72").Formula = "=H72-I72"
$ws.Range("K72").Formula = "=B72+E72+H72"
$ws.Range("L72").Formula = "=C72+F72+I72"
$ws.Range("M72").Formula = "=K72-L72"

$ws.Range("A73").Value = 43967
$ws.Range("B73").Value = 70
$ws.Range("C73").Value = 25
$ws.Range("D73").Formula = "=B73-C73"
$ws.Range("E73").Value = 160
$ws.Range("F73").Value = 38
$ws.Range("G73").Formula = "=E73-F73"
$ws.Range("H73").Value = 20
$ws.Range("I73").Value = 8
$ws.Range("J73").Formula = "=H73-I73"
$ws.Range("K73").Formula = "=B73+E73+H73"
$ws.Range("L73").Formula = "=C73+F73+I73"
$ws.Range("M73").Formula = "=K73-L73"

$ws.Range("A74").Value = 43968
$ws.Range("B74").Value = 67
$ws.Range("C74").Value = 25
$ws.Range("D74").Formula = "=B74-C74"
$ws.Range("E74").Value = 149
$ws.Range("F74").Value = 42
$ws.Range("G74").Formula = "=E74-F74"
$ws.Range("H74").Value = 21
$ws.Range("I74").Value = 8
$ws.Range("J74").Formula = "=H74-I74"
$ws.Range("K74").Formula = "=B74+E74+H74"
$ws.Range("L74").Formula = "=C74+F74+I74"
$ws.Range("M74").Formula = "=K74-L74"

$ws.Range("A75").Value = 43969
$ws.Range("B75").Value = 63
$ws.Range("C75").Value = 23
$ws.Range("D75").Formula = "=B75-C75"
$ws.Range("E75").Value = 147
$ws.Range("F75").Value = 42
$ws.Range("G75").Formula = "=E75-F75"
$ws.Range("H75").Value = 20
$ws.Range("I75").Value = 8
$ws.Range("J75").Formula = "=H75-I75"
$ws.Range("K75").Formula = "=B75+E75+H75"
$ws.Range("L75").Formula = "=C75+F75+I75"
$ws.Range("M75").Formula = "=K75-L75"

$ws.Range("A76").Value = 43970
$ws.Range("B76").Value = 59
$ws.Range("C76").Value = 23
$ws.Range("D76").Formula = "=B76-C76"
$ws.Range("E76").Value = 141
$ws.Range("F76").Value = 35
$ws.Range("G76").Formula = "=E76-F76"
$ws.Range("H76").Value = 17
$ws.Range("I76").Value = 4
$ws.Range("J76").Formula = "=H76-I76"
$ws.Range("K76").Formula = "=B76+E76+H76"
$ws.Range("L76").Formula = "=C76+F76+I76"
$ws.Range("M76").Formula = "=K76-L76"

# --- Scroll the view down to the new data and select the last entered cell ---
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J76").Select()
